$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: round coordinates to whole numbers, clear Starttid/Sluttid ---
$ws.Range("Q2").Value = 469609
$ws.Range("R2").Value = 7039805
$ws.Range("Z2").Value = ""
$ws.Range("AB2").Value = ""

# --- Rows 3 & 4: swap species data (A,B,D,E,F,G,H) and update coordinates ---
$a3 = $ws.Range("A3").Value()
$b3 = $ws.Range("B3").Value()
$d3 = $ws.Range("D3").Value()
$e3 = $ws.Range("E3").Value()
$f3 = $ws.Range("F3").Value()
$g3 = $ws.Range("G3").Value()
$h3 = $ws.Range("H3").Value()

$a4 = $ws.Range("A4").Value()
$b4 = $ws.Range("B4").Value()
$d4 = $ws.Range("D4").Value()
$e4 = $ws.Range("E4").Value()
$f4 = $ws.Range("F4").Value()
$g4 = $ws.Range("G4").Value()
$h4 = $ws.Range("H4").Value()

$ws.Range("A3").Value = $a4
$ws.Range("B3").Value = $b4
$ws.Range("D3").Value = $d4
$ws.Range("E3").Value = $e4
$ws.Range("F3").Value = $f4
$ws.Range("G3").Value = $g4
$ws.Range("H3").Value = $h4

$ws.Range("A4").Value = $a3
$ws.Range("B4").Value = $b3
$ws.Range("D4").Value = $d3
$ws.Range("E4").Value = $e3
$ws.Range("F4").Value = $f3
$ws.Range("G4").Value = $g3
$ws.Range("H4").Value = $h3

# Coordinates: row3 gets rounded row4-original, row4 gets rounded row3-original
$ws.Range("Q3").Value = 469608
$ws.Range("R3").Value = 7039809
$ws.Range("Q4").Value = 469597
$ws.Range("R4").Value = 7039829

# Clear Starttid/Sluttid for rows 3 and 4
$ws.Range("Z3").Value = ""
$ws.Range("AB3").Value = ""
$ws.Range("Z4").Value = ""
$ws.Range("AB4").Value = ""
